$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9496732950210571
$ws.Range("B1").Value = 1.710909962654114
$ws.Range("C1").Value = 4.634339332580566
$ws.Range("D1").Value = 1.532417058944702
$ws.Range("E1").Value = 1.097716093063354
